# Update "想去人数" (F column) values on the 展览 (rId1) and 全部类型 (rId4) sheets
# to reflect freshly generated data (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 13850
$wsExhibit.Range("F4").Value = 13627
$wsExhibit.Range("F19").Value = 537
$wsExhibit.Range("F21").Value = 415
$wsExhibit.Range("F24").Value = 842
$wsExhibit.Range("F25").Value = 98

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 13850
$wsAll.Range("F5").Value = 13627
$wsAll.Range("F26").Value = 537
$wsAll.Range("F28").Value = 415
$wsAll.Range("F31").Value = 842
$wsAll.Range("F37").Value = 98
